# Update Koffi benchmark results and instructions
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Linux
$ws2 = $wb.Worksheets.Item(2)   # Windows

# --- Linux sheet benchmark values (column C drives the "x N" labels in D via formula) ---
$ws1.Range("C4").Value = 0.61
$ws1.Range("C8").Value = 0.59
$ws1.Range("C9").Value = 0.008
$ws1.Range("C12").Value = 0.94
$ws1.Range("C13").Value = 0.3

# --- Windows sheet benchmark values ---
$ws2.Range("C4").Value = 0.64
$ws2.Range("C8").Value = 0.55
$ws2.Range("C9").Value = 0.01
$ws2.Range("C12").Value = 0.92
$ws2.Range("C13").Value = 0.28

# Recalculate formulas/charts so everything derived from the cells above is current.
$excel.CalculateFullRebuild()
$wb.RefreshAll()

# --- Update view/selection state to match the author's final position ---
# The Linux sheet is scrolled/selected first, but it is not the tab left active.
$ws1.Activate()
$ws1.Range("O20").Select()

# The Windows sheet ends up being the active / selected tab, with C14 selected.
$ws2.Activate()
$ws2.Range("C14").Select()
